$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 19.434592
$ws.Range("H2").Value = 58.303776
$ws.Range("I2").Value = 0.1244167820899015
$ws.Range("J2").Value = 0.1244167820899015
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 40.92730033333333
$ws.Range("N2").Value = 122.781901
$ws.Range("O2").Value = 0.3921621516522625
$ws.Range("P2").Value = 0.3921621516522625
$ws.Range("Q2").Value = 795.4053836397972
$ws.Range("R2").Value = 7158.648452758175
$ws.Range("S2").Value = 0.04879155296602647
$ws.Range("T2").Value = 0.04879155296602646

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 19.434592
$ws.Range("H3").Value = 58.303776
$ws.Range("I3").Value = 0.1244167820899015
$ws.Range("J3").Value = 0.1244167820899015
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 7.598228333333334
$ws.Range("N3").Value = 22.794685
$ws.Range("O3").Value = 0.07280562235174674
$ws.Range("P3").Value = 0.07280562235174673
$ws.Range("Q3").Value = 147.6684675811733
$ws.Range("R3").Value = 1329.01620823056
$ws.Range("S3").Value = 0.009058241251056939
$ws.Range("T3").Value = 0.009058241251056936

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 19.434592
$ws.Range("H4").Value = 58.303776
$ws.Range("I4").Value = 0.1244167820899015
$ws.Range("J4").Value = 0.1244167820899015
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 49.40125466666666
$ws.Range("N4").Value = 148.203764
$ws.Range("O4").Value = 0.4733589112063359
$ws.Range("P4").Value = 0.4733589112063358
$ws.Range("Q4").Value = 960.0932287347624
$ws.Range("R4").Value = 8640.839058612863
$ws.Range("S4").Value = 0.05889379250587174
$ws.Range("T4").Value = 0.05889379250587172

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 19.434592
$ws.Range("H5").Value = 58.303776
$ws.Range("I5").Value = 0.1244167820899015
$ws.Range("J5").Value = 0.1244167820899015
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 6.436425
$ws.Range("N5").Value = 19.309275
$ws.Range("O5").Value = 0.06167331478965488
$ws.Range("P5").Value = 0.06167331478965488
$ws.Range("Q5").Value = 125.0892938136
$ws.Range("R5").Value = 1125.8036443224
$ws.Range("S5").Value = 0.007673195366946393
$ws.Range("T5").Value = 0.007673195366946392

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 48.891945
$ws.Range("H6").Value = 146.675835
$ws.Range("I6").Value = 0.3129974875220664
$ws.Range("J6").Value = 0.3129974875220664
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 40.92730033333333
$ws.Range("N6").Value = 122.781901
$ws.Range("O6").Value = 0.3921621516522625
$ws.Range("P6").Value = 0.3921621516522625
$ws.Range("Q6").Value = 2001.015316895815
$ws.Range("R6").Value = 18009.13785206233
$ws.Range("S6").Value = 0.1227457681684058
$ws.Range("T6").Value = 0.1227457681684057

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 48.891945
$ws.Range("H7").Value = 146.675835
$ws.Range("I7").Value = 0.3129974875220664
$ws.Range("J7").Value = 0.3129974875220664
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 7.598228333333334
$ws.Range("N7").Value = 22.794685
$ws.Range("O7").Value = 0.07280562235174674
$ws.Range("P7").Value = 0.07280562235174673
$ws.Range("Q7").Value = 371.4921617707751
$ws.Range("R7").Value = 3343.429455936976
$ws.Range("S7").Value = 0.02278797687357713
$ws.Range("T7").Value = 0.02278797687357712

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 48.891945
$ws.Range("H8").Value = 146.675835
$ws.Range("I8").Value = 0.3129974875220664
$ws.Range("J8").Value = 0.3129974875220664
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 49.40125466666666
$ws.Range("N8").Value = 148.203764
$ws.Range("O8").Value = 0.4733589112063359
$ws.Range("P8").Value = 0.4733589112063358
$ws.Range("Q8").Value = 2415.32342609366
$ws.Range("R8").Value = 21737.91083484294
$ws.Range("S8").Value = 0.1481601499037641
$ws.Range("T8").Value = 0.148160149903764

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 48.891945
$ws.Range("H9").Value = 146.675835
$ws.Range("I9").Value = 0.3129974875220664
$ws.Range("J9").Value = 0.3129974875220664
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 6.436425
$ws.Range("N9").Value = 19.309275
$ws.Range("O9").Value = 0.06167331478965488
$ws.Range("P9").Value = 0.06167331478965488
$ws.Range("Q9").Value = 314.689337096625
$ws.Range("R9").Value = 2832.204033869625
$ws.Range("S9").Value = 0.01930359257631948
$ws.Range("T9").Value = 0.01930359257631947

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 65.19353
$ws.Range("H10").Value = 195.58059
$ws.Range("I10").Value = 0.4173573191390618
$ws.Range("J10").Value = 0.4173573191390618
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 40.92730033333333
$ws.Range("N10").Value = 122.781901
$ws.Range("O10").Value = 0.3921621516522625
$ws.Range("P10").Value = 0.3921621516522625
$ws.Range("Q10").Value = 2668.195182100176
$ws.Range("R10").Value = 24013.75663890159
$ws.Range("S10").Value = 0.1636717442813945
$ws.Range("T10").Value = 0.1636717442813945

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 65.19353
$ws.Range("H11").Value = 195.58059
$ws.Range("I11").Value = 0.4173573191390618
$ws.Range("J11").Value = 0.4173573191390618
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 7.598228333333334
$ws.Range("N11").Value = 22.794685
$ws.Range("O11").Value = 0.07280562235174674
$ws.Range("P11").Value = 0.07280562235174673
$ws.Range("Q11").Value = 495.3553267960167
$ws.Range("R11").Value = 4458.197941164151
$ws.Range("S11").Value = 0.03038595936297598
$ws.Range("T11").Value = 0.03038595936297597

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 65.19353
$ws.Range("H12").Value = 195.58059
$ws.Range("I12").Value = 0.4173573191390618
$ws.Range("J12").Value = 0.4173573191390618
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 49.40125466666666
$ws.Range("N12").Value = 148.203764
$ws.Range("O12").Value = 0.4733589112063359
$ws.Range("P12").Value = 0.4733589112063358
$ws.Range("Q12").Value = 3220.642178148973
$ws.Range("R12").Value = 28985.77960334076
$ws.Range("S12").Value = 0.1975598061716616
$ws.Range("T12").Value = 0.1975598061716615

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 65.19353
$ws.Range("H13").Value = 195.58059
$ws.Range("I13").Value = 0.4173573191390618
$ws.Range("J13").Value = 0.4173573191390618
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 6.436425
$ws.Range("N13").Value = 19.309275
$ws.Range("O13").Value = 0.06167331478965488
$ws.Range("P13").Value = 0.06167331478965488
$ws.Range("Q13").Value = 419.61326633025
$ws.Range("R13").Value = 3776.51939697225
$ws.Range("S13").Value = 0.02573980932302981
$ws.Range("T13").Value = 0.02573980932302981

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 22.685484
$ws.Range("H14").Value = 68.05645200000001
$ws.Range("I14").Value = 0.1452284112489703
$ws.Range("J14").Value = 0.1452284112489703
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 40.92730033333333
$ws.Range("N14").Value = 122.781901
$ws.Range("O14").Value = 0.3921621516522625
$ws.Range("P14").Value = 0.3921621516522625
$ws.Range("Q14").Value = 928.455616875028
$ws.Range("R14").Value = 8356.100551875252
$ws.Range("S14").Value = 0.05695308623643584
$ws.Range("T14").Value = 0.05695308623643583

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 22.685484
$ws.Range("H15").Value = 68.05645200000001
$ws.Range("I15").Value = 0.1452284112489703
$ws.Range("J15").Value = 0.1452284112489703
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 7.598228333333334
$ws.Range("N15").Value = 22.794685
$ws.Range("O15").Value = 0.07280562235174674
$ws.Range("P15").Value = 0.07280562235174673
$ws.Range("Q15").Value = 172.36948728418
$ws.Range("R15").Value = 1551.32538555762
$ws.Range("S15").Value = 0.0105734448641367
$ws.Range("T15").Value = 0.0105734448641367

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 22.685484
$ws.Range("H16").Value = 68.05645200000001
$ws.Range("I16").Value = 0.1452284112489703
$ws.Range("J16").Value = 0.1452284112489703
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 49.40125466666666
$ws.Range("N16").Value = 148.203764
$ws.Range("O16").Value = 0.4733589112063359
$ws.Range("P16").Value = 0.4733589112063358
$ws.Range("Q16").Value = 1120.691372320592
$ws.Range("R16").Value = 10086.22235088533
$ws.Range("S16").Value = 0.06874516262503856
$ws.Range("T16").Value = 0.06874516262503855

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 22.685484
$ws.Range("H17").Value = 68.05645200000001
$ws.Range("I17").Value = 0.1452284112489703
$ws.Range("J17").Value = 0.1452284112489703
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 6.436425
$ws.Range("N17").Value = 19.309275
$ws.Range("O17").Value = 0.06167331478965488
$ws.Range("P17").Value = 0.06167331478965488
$ws.Range("Q17").Value = 146.0134163547
$ws.Range("R17").Value = 1314.1207471923
$ws.Range("S17").Value = 0.008956717523359201
$ws.Range("T17").Value = 0.008956717523359199
